$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed as a number by
# Excel's type inference (e.g. "1.004", "41.37") need NumberFormat = "@" applied
# first so the literal text (including trailing zeros) survives, matching the
# scraped inline-string cell contents exactly.
$forceTextCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D12", "D14", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Rows 44-45 swapped rank order (EnergySwap now ranks above Decentraland),
# each also carrying refreshed Price / Volume(1h) figures.
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "13.23"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5988"
$ws.Range("E45").Value = "  +1.53%  "

# Refreshed Price (D) / Volume(1h) (E) figures for every other coin row.
$ws.Range("D2").Value = "28.152.53"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.800.87"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "313.80"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.5266"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").Value = "0.3812"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "0.07959"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").Value = "41.37"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "6.326"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "20.60"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "1.807.91"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "7.327"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "92.69"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "0.00001091"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Value = "0.06614"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "5.971"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "28.187.40"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "2.238"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "158.34"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "20.53"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").Value = "2.407"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "2.008.72"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "122.95"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "0.1094"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "1.057"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").Value = "3.666"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "5.534"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("D35").Value = "0.07261"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").Value = "12.19"
$ws.Range("E36").Value = "  +8.65%  "
$ws.Range("D37").Value = "8.879"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "0.2165"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "0.02307"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "5.048"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "0.6183"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "1.165"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "1.370"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D46").Value = "3.762"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "126.42"
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").Value = "1.201"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "1.928"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "0.06824"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "72.77"
$ws.Range("E51").Value = "  -2.18%  "
